# Auto-generated script applying scheduled market-data refresh updates
# to the Diabolos_Profits workbook, sheet by sheet.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3347.524
$ws.Range("J40").Value = 3547
$ws.Range("L40").Value = 3547
$ws.Range("N40").Value = -3897
$ws.Range("H62").Value = 25999.44
$ws.Range("I62").Value = 5587.4116
$ws.Range("K62").Value = 5587.4116
$ws.Range("M62").Value = -4963.4116
$ws.Range("H65").Value = 25999.44
$ws.Range("I65").Value = 5587.4116
$ws.Range("K65").Value = 27937.058
$ws.Range("M65").Value = -24817.058
$ws.Range("H137").Value = 4224.839
$ws.Range("I137").Value = 3118.6
$ws.Range("J137").Value = 6236.1816
$ws.Range("K137").Value = 9355.799999999999
$ws.Range("L137").Value = 18708.5448
$ws.Range("M137").Value = -6805.799999999999
$ws.Range("N137").Value = -23808.5448
$ws.Range("H138").Value = 3487.64
$ws.Range("J138").Value = 3450.6052
$ws.Range("L138").Value = 10351.8156
$ws.Range("N138").Value = -20631.8156

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16036.289
$ws.Range("I32").Value = 8714.809999999999
$ws.Range("J32").Value = 39098.95
$ws.Range("K32").Value = 8714.809999999999
$ws.Range("L32").Value = 39098.95
$ws.Range("M32").Value = -8427.809999999999
$ws.Range("N32").Value = -39672.95
$ws.Range("H61").Value = 3161.3928
$ws.Range("I61").Value = 2105.1667
$ws.Range("J61").Value = 9498.75
$ws.Range("K61").Value = 2105.1667
$ws.Range("L61").Value = 9498.75
$ws.Range("M61").Value = -1893.1667
$ws.Range("N61").Value = -9922.75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("H132").Value = 2692.9285
$ws.Range("I132").Value = 2291.162
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 6873.485999999999
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -4343.485999999999
$ws.Range("N132").Value = -22058
$ws.Range("H136").Value = 3161.3928
$ws.Range("I136").Value = 2105.1667
$ws.Range("J136").Value = 9498.75
$ws.Range("K136").Value = 6315.500100000001
$ws.Range("L136").Value = 28496.25
$ws.Range("M136").Value = -3765.500100000001
$ws.Range("N136").Value = -33596.25
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 60849.8
$ws.Range("J137").Value = 60849.8
$ws.Range("L137").Value = 60849.8
$ws.Range("N137").Value = -71049.8

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 792.6111
$ws.Range("I22").Value = 387.9091
$ws.Range("J22").Value = 1428.5714
$ws.Range("K22").Value = 387.9091
$ws.Range("L22").Value = 1428.5714
$ws.Range("M22").Value = -37.90910000000002
$ws.Range("N22").Value = -2128.5714
$ws.Range("H31").Value = 4092.1
$ws.Range("I31").Value = 3469.8
$ws.Range("J31").Value = 4216.56
$ws.Range("K31").Value = 3469.8
$ws.Range("L31").Value = 4216.56
$ws.Range("M31").Value = -3174.8
$ws.Range("N31").Value = -4806.56
$ws.Range("H34").Value = 4092.1
$ws.Range("I34").Value = 3469.8
$ws.Range("J34").Value = 4216.56
$ws.Range("K34").Value = 3469.8
$ws.Range("L34").Value = 4216.56
$ws.Range("M34").Value = -3267.8
$ws.Range("N34").Value = -4620.56
$ws.Range("H41").Value = 29020
$ws.Range("J41").Value = 33530
$ws.Range("L41").Value = 33530
$ws.Range("N41").Value = -34386
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("H58").Value = 315503.7
$ws.Range("I58").Value = 2086.611
$ws.Range("J58").Value = 718468.5
$ws.Range("K58").Value = 2086.611
$ws.Range("L58").Value = 718468.5
$ws.Range("M58").Value = -1883.611
$ws.Range("N58").Value = -718874.5
$ws.Range("H93").Value = 11108.1
$ws.Range("I93").Value = 11108.1
$ws.Range("K93").Value = 11108.1
$ws.Range("M93").Value = -9236.1
$ws.Range("H122").Value = 2806.1667
$ws.Range("I122").Value = 2986.4
$ws.Range("J122").Value = 1905
$ws.Range("K122").Value = 8959.200000000001
$ws.Range("L122").Value = 5715
$ws.Range("M122").Value = -6509.200000000001
$ws.Range("N122").Value = -10615
$ws.Range("H132").Value = 210024.73
$ws.Range("I132").Value = 1268.4722
$ws.Range("K132").Value = 3805.4166
$ws.Range("M132").Value = -1275.4166
$ws.Range("H134").Value = 3655.2454
$ws.Range("I134").Value = 3238.7297
$ws.Range("K134").Value = 9716.1891
$ws.Range("M134").Value = -7181.1891
$ws.Range("H136").Value = 315503.7
$ws.Range("I136").Value = 2086.611
$ws.Range("J136").Value = 718468.5
$ws.Range("K136").Value = 6259.833
$ws.Range("L136").Value = 2155405.5
$ws.Range("M136").Value = -3709.833
$ws.Range("N136").Value = -2160505.5
$ws.Range("H139").Value = 84854.5
$ws.Range("J139").Value = 84854.5
$ws.Range("L139").Value = 84854.5
$ws.Range("N139").Value = -95134.5
$ws.Range("L52").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 886.0714
$ws.Range("J2").Value = 1125.5
$ws.Range("L2").Value = 6753
$ws.Range("N2").Value = -6979
$ws.Range("H80").Value = 6130.5654
$ws.Range("J80").Value = 6375.375
$ws.Range("L80").Value = 19126.125
$ws.Range("N80").Value = -20998.125
$ws.Range("H83").Value = 6130.5654
$ws.Range("J83").Value = 6375.375
$ws.Range("L83").Value = 57378.375
$ws.Range("N83").Value = -66738.375
$ws.Range("H92").Value = 555.4
$ws.Range("I92").Value = 530.6667
$ws.Range("K92").Value = 1592.0001
$ws.Range("M92").Value = -344.0001
$ws.Range("H98").Value = 649.5
$ws.Range("J98").Value = 649.5
$ws.Range("L98").Value = 1948.5
$ws.Range("N98").Value = -4944.5
$ws.Range("H107").Value = 236.64706
$ws.Range("I107").Value = 166.66667
$ws.Range("J107").Value = 251.64285
$ws.Range("K107").Value = 500.00001
$ws.Range("L107").Value = 754.9285500000001
$ws.Range("M107").Value = 1419.99999
$ws.Range("N107").Value = -4594.928550000001
$ws.Range("H131").Value = 11028.32
$ws.Range("I131").Value = 2236.2856
$ws.Range("K131").Value = 6708.8568
$ws.Range("M131").Value = -1668.8568
$ws.Range("H132").Value = 2518.2458
$ws.Range("I132").Value = 1919.1818
$ws.Range("J132").Value = 2650.04
$ws.Range("K132").Value = 17272.6362
$ws.Range("L132").Value = 23850.36
$ws.Range("M132").Value = -14742.6362
$ws.Range("N132").Value = -28910.36

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2102.718
$ws.Range("I102").Value = 992.2083
$ws.Range("K102").Value = 992.2083
$ws.Range("M102").Value = 629.7917
$ws.Range("H122").Value = 152204.66
$ws.Range("I122").Value = 186956.56
$ws.Range("K122").Value = 560869.6799999999
$ws.Range("M122").Value = -558419.6799999999
$ws.Range("H126").Value = 6981.552
$ws.Range("I126").Value = 8759.166999999999
$ws.Range("J126").Value = 4072.7273
$ws.Range("K126").Value = 26277.501
$ws.Range("L126").Value = 12218.1819
$ws.Range("M126").Value = -23807.501
$ws.Range("N126").Value = -17158.1819
$ws.Range("H140").Value = 171035.8
$ws.Range("I140").Value = 103472.664
$ws.Range("J140").Value = 199991.42
$ws.Range("K140").Value = 103472.664
$ws.Range("L140").Value = 199991.42
$ws.Range("M140").Value = -98292.664
$ws.Range("N140").Value = -210351.42

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16954082
$ws.Range("I132").Value = 25644272
$ws.Range("K132").Value = 76932816
$ws.Range("M132").Value = -76930286
$ws.Range("H136").Value = 6193.224
$ws.Range("I136").Value = 6188.6924
$ws.Range("K136").Value = 18566.0772
$ws.Range("M136").Value = -16016.0772

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 6333.3335
$ws.Range("J5").Value = 6333.3335
$ws.Range("L5").Value = 6333.3335
$ws.Range("N5").Value = -6557.3335
$ws.Range("H93").Value = 46884.332
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -64992
$ws.Range("H122").Value = 1817.7273
$ws.Range("I122").Value = 1493.3334
$ws.Range("K122").Value = 4480.0002
$ws.Range("M122").Value = -2030.0002
$ws.Range("H132").Value = 569027.75
$ws.Range("I132").Value = 805350.7
$ws.Range("J132").Value = 31930.182
$ws.Range("K132").Value = 2416052.1
$ws.Range("L132").Value = 95790.546
$ws.Range("M132").Value = -2413522.1
$ws.Range("N132").Value = -100850.546
